# Auto-generated Excel COM-interop script to apply numeric updates
# to the "北京-漫展信息" workbook per the provided diff.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 657
$ws.Range("F7").Value = 2172
$ws.Range("F8").Value = 905
$ws.Range("F9").Value = 864
$ws.Range("F10").Value = 417
$ws.Range("F11").Value = 98
$ws.Range("F13").Value = 326
$ws.Range("F15").Value = 1105
$ws.Range("F18").Value = 1827
$ws.Range("F19").Value = 46
$ws.Range("F20").Value = 31
$ws.Range("F21").Value = 27
$ws.Range("F24").Value = 1463
$ws.Range("F25").Value = 12
$ws.Range("F26").Value = 534
$ws.Range("F28").Value = 634
$ws.Range("F29").Value = 428
$ws.Range("F30").Value = 2536
$ws.Range("G30").Value = 72
$ws.Range("F31").Value = 389
$ws.Range("F32").Value = 99
$ws.Range("F33").Value = 1402
$ws.Range("F34").Value = 609
$ws.Range("F35").Value = 486
$ws.Range("F36").Value = 205
$ws.Range("F37").Value = 938
$ws.Range("F38").Value = 713
$ws.Range("F40").Value = 535
$ws.Range("F41").Value = 530

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F14").Value = 76
$ws.Range("F22").Value = 127

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 230
$ws.Range("F3").Value = 2925
$ws.Range("F6").Value = 334

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 230
$ws.Range("F8").Value = 657
$ws.Range("F11").Value = 905
$ws.Range("F12").Value = 864
$ws.Range("F13").Value = 417
$ws.Range("F14").Value = 98
$ws.Range("F16").Value = 326
$ws.Range("F19").Value = 1105
$ws.Range("F23").Value = 334
$ws.Range("F24").Value = 1827
$ws.Range("F25").Value = 46
$ws.Range("F31").Value = 76
$ws.Range("F32").Value = 1463
$ws.Range("F34").Value = 12
$ws.Range("F35").Value = 534
$ws.Range("F37").Value = 428
$ws.Range("F38").Value = 2536
$ws.Range("G38").Value = 72
$ws.Range("F39").Value = 99
$ws.Range("F40").Value = 609
$ws.Range("F41").Value = 486
$ws.Range("F42").Value = 205
$ws.Range("F43").Value = 938
$ws.Range("F46").Value = 713
$ws.Range("F48").Value = 535
$ws.Range("F49").Value = 531
